$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with new word data
$ws.Range("B2").Value = "meandering"
$ws.Range("C2").Value = "To wind or turn in a course or passage; to be intricate.;To wind, turn, or twist; to make flexuous."
$ws.Range("D2").Value = "блуждающий"

# Remove row 3 entirely (the "faint" entry)
$ws.Rows("3:3").Delete()
